$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rowIndex = 5

$ws.Cells.Item($rowIndex, 1).Value = 42608.893449074072
$ws.Cells.Item($rowIndex, 1).NumberFormat = "m/d/yy h:mm"

$ws.Cells.Item($rowIndex, 2).Value = 24
$ws.Cells.Item($rowIndex, 3).Value = 58
$ws.Cells.Item($rowIndex, 4).Value = 38
$ws.Cells.Item($rowIndex, 5).Value = 99
$ws.Cells.Item($rowIndex, 6).Value = 0
$ws.Cells.Item($rowIndex, 7).Value = 22588
$ws.Cells.Item($rowIndex, 8).Value = 18672
$ws.Cells.Item($rowIndex, 9).Value = 1173
$ws.Cells.Item($rowIndex, 10).Value = 190
$ws.Cells.Item($rowIndex, 11).Value = 124
$ws.Cells.Item($rowIndex, 12).Value = 3
$ws.Cells.Item($rowIndex, 13).Value = 0
$ws.Cells.Item($rowIndex, 14).Value = "Named"
